# [FEATURE] working on import
#
# Adds a "Resources" sheet after "Meta", repopulates "Meta" with the
# Sheet / Assembly-Qualified-Type-Name mapping row, and fills "Resources"
# with the Name / TwoLetterISOLanguageName / Value rows.

$wb = $excel.ActiveWorkbook

# --- Recreate "Meta" (delete + re-add) -------------------------------
# This bumps the sheet's internal sheetId the same way the original
# authoring session did (Meta: sheetId 1 -> 2, new Resources: sheetId 3).
$oldMeta = $wb.Sheets("Meta")
$newMeta = $wb.Worksheets.Add($null, $oldMeta)
$newMeta.Name = "MetaRebuilt"
$oldMeta.Delete()
$wb.Sheets("MetaRebuilt").Name = "Meta"

$meta = $wb.Sheets("Meta")

# Header row (bold)
$meta.Range("A1").Value = "Sheet"
$meta.Range("B1").Value = "Assembly Qualified Type Name"
$meta.Range("A1:B1").Font.Bold = $true

# Data row
$meta.Range("A2").Value = "Resources"
$meta.Range("B2").Value = "Hydra.Nh.Infrastructure.I18n.ResourceItem, Hydra.Nh"

# Column B width (character-width units; the saved file's column width
# ends up 11/12 wider than the value assigned here).
$meta.Columns("B").ColumnWidth = 57.083333333333336

# --- Add "Resources" sheet after "Meta" -------------------------------
$resources = $wb.Worksheets.Add($null, $wb.Sheets("Meta"))
$resources.Name = "Resources"

$resources.Range("A1").Value = "Name"
$resources.Range("B1").Value = "TwoLetterISOLanguageName"
$resources.Range("C1").Value = "Value"

$resources.Range("A2").Value = "FOO"
$resources.Range("B2").Value = "DE"
$resources.Range("C2").Value = "FOO_DE"

$resources.Range("A3").Value = "FOO"
$resources.Range("B3").Value = "EN"
$resources.Range("C3").Value = "FOO_EN"

$resources.Range("A4").Value = "BAR.BAZ"
$resources.Range("B4").Value = "DE"
$resources.Range("C4").Value = "BAR.BAZ_DE"

$resources.Range("A5").Value = "BAR.BAZ"
$resources.Range("B5").Value = "EN"
$resources.Range("C5").Value = "BAR.BAZ_EN"

$resources.Columns("B").ColumnWidth = 26.083333333333332

# --- Selections / active sheet ----------------------------------------
$meta.Range("B7").Select()
$resources.Range("C6").Select()
